$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B, C, D, E in this sheet hold text data (coin name, link, price, volume%).
# Force text format on column D (price) so numeric-looking strings like "242.41"
# are not auto-converted to numbers by Excel, matching the original inline-string cells.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.386.02"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.883.26"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.41"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.08033"
$ws.Range("E8").Value = "  +3.41%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.889.36"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.248"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7201"
$ws.Range("E14").Value = "  +1.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.37"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.335"
$ws.Range("E16").Value = "  +5.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008494"
$ws.Range("E17").Value = "  +3.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.396.55"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.58"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.149.28"
$ws.Range("E20").Value = "  +0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.26"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +0.16%  "
$ws.Range("E23").Value = "  +0.41%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("E25").Value = "  -2.31%  "
$ws.Range("E26").Value = "  +0.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.060"
$ws.Range("E27").Value = "  +0.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.58"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.510"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.419"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.344"
$ws.Range("E31").Value = "  +0.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.199"
$ws.Range("E32").Value = "  -6.27%  "
$ws.Range("E33").Value = "  +2.55%  "
$ws.Range("E34").Value = "  +1.02%  "
$ws.Range("E35").Value = "  +0.56%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7504"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.706"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.290.69"
$ws.Range("E38").Value = "  +9.84%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01889"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.744"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.606"
$ws.Range("E41").Value = "  +3.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9289"
$ws.Range("E42").Value = "  +4.48%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "111.94"
$ws.Range("E43").Value = "  +5.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "74.05"
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("E45").Value = "  +0.10%  "
$ws.Range("E46").Value = "  +5.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.036.68"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.807"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.5221"
$ws.Range("E49").Value = "  +0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.504"
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4397"
$ws.Range("E51").Value = "  +2.11%  "
